# Updates cryptos list values (price + volume/1h columns) to match the
# latest scrape. Rows 39/40 (Hedera / TheSandbox) also swap places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text updates (values that Excel would never mistake for numbers) ---
$ws.Range("D2").Value = "28.037.91"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "1.913.99"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "1.889.85"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").Value = "28.063.09"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").Value = "2.148.66"
$ws.Range("E23").Value = "  +2.79%  "
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("E26").Value = "  -1.58%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  +5.34%  "
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("E36").Value = "  -4.25%  "
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  +5.31%  "
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("E51").Value = "  -0.04%  "

# --- Numeric-looking text updates ---
# These "Price" cells look like plain decimals (e.g. "315.65"), so a naive
# Value assignment would have Excel reinterpret them as actual numbers and
# mangle formatting (trailing zeros dropped, floating point noise, etc).
# Force the cell to Text first, assign, then drop back to the default style
# so the cell keeps looking exactly like its neighbours (no leftover quote-
# prefix / custom number format on the cell).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.65"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4812"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3810"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07362"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9343"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.82"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07806"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.500"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.636"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.97"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008864"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.76"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.169"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.63"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.910"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.50"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.130"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.68"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.968"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08950"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.298"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.255"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7739"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.674"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02050"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5525"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05305"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.022"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1526"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.498"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.70"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "108.56"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4828"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.647"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.00"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06084"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
